# Update Name of Algo
# Applies the updated KNN imputation results to columns B and D of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -7.547
$ws.Range("B8").Value = 6.467000000000001
$ws.Range("B10").Value = 6.483000000000001
$ws.Range("D11").Value = -7.342000000000001
$ws.Range("B12").Value = 5.553
$ws.Range("D12").Value = -6.985000000000001
$ws.Range("D15").Value = -8.317
$ws.Range("D17").Value = -8.126999999999999
$ws.Range("B18").Value = 5.137
$ws.Range("B25").Value = 5.595999999999999
$ws.Range("D26").Value = -6.917
$ws.Range("D27").Value = -7.924000000000001
$ws.Range("D28").Value = -8.002000000000001
$ws.Range("D32").Value = -6.781000000000001
$ws.Range("B37").Value = 8.626999999999999
$ws.Range("D37").Value = -7.967000000000001
$ws.Range("D41").Value = -8.061000000000002
$ws.Range("D47").Value = -7.664
$ws.Range("D51").Value = -8.372999999999999
$ws.Range("B55").Value = 4.722
$ws.Range("D65").Value = -7.704000000000001
$ws.Range("B68").Value = 5.084999999999999
$ws.Range("D73").Value = -7.874000000000001
$ws.Range("B77").Value = 5.186999999999999
$ws.Range("B78").Value = 7.472
$ws.Range("B79").Value = 5.077
$ws.Range("B80").Value = 7.343999999999999
$ws.Range("B81").Value = 6.439
$ws.Range("B82").Value = 5.464
$ws.Range("B84").Value = 6.571
$ws.Range("D84").Value = -8.228000000000002
$ws.Range("D85").Value = -8.740999999999998
$ws.Range("D89").Value = -8.068999999999999
$ws.Range("D93").Value = -6.873
$ws.Range("D95").Value = -7.571000000000001
$ws.Range("D98").Value = -7.027000000000001
$ws.Range("D99").Value = -8.273999999999999
$ws.Range("B101").Value = 6.074000000000001
$ws.Range("D101").Value = -7.991
$ws.Range("B102").Value = 7.313999999999998
$ws.Range("D102").Value = -7.683999999999999
